$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update Sheet1's last two data rows (still rows 4 & 5, before the header row gets
#     inserted below) ---
#     Row4: A "EugenBorisik"   -> "EugenBorisik2"   (B stays "qwerty12345000")
#     Row5: A "NataliaDamorad" -> "NataliaDamorad2" (B stays "Vintage20")
$ws1.Range("A4").Value = "EugenBorisik2"
$ws1.Range("A5").Value = "NataliaDamorad2"

# --- Duplicate Sheet1 into a new Sheet2 placed right after it (this carries over the
#     current data/formatting, which we then touch up independently) ---
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"
$ws2.Columns.Item(1).ColumnWidth = 17.85546875
$ws2.Columns.Item(2).ColumnWidth = 21.85546875

# --- Insert the "UserName"/"Password" header row at the top of both sheets ---
$ws1.Rows.Item(1).Insert()
$ws1.Range("A1").Value = "UserName"
$ws1.Range("B1").Value = "Password"

$ws2.Rows.Item(1).Insert()
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"

# --- Selections: whole data range on each sheet ---
$ws1.Range("A1:B6").Select()
$ws2.Range("A1:B6").Select()

# --- Turn both ranges into native Excel tables ---
$tbl1 = $ws1.ListObjects.Add(1, $ws1.Range("A1:B6"), $null, 1)
$tbl1.Name = "Information"

$tbl2 = $ws2.ListObjects.Add(1, $ws2.Range("A1:B6"), $null, 1)
$tbl2.Name = "Information2"

# Sheet2 is the tab that ends up active/selected, matching the source edit.
$ws2.Select()
